# "new future: gabari summary" -- extend the "خروجی" (output) sheet with
# a full multi-leg gabari (loading-gauge) summary table: fix up the header
# wording, drop the now-unused "F2" free-space entry-size cell on the first
# leg, add a 7th column for the new "فضای سازه" (structure space) check,
# and append four more route legs (rows 3-6), each formatted the same way
# as the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row (row 1) -----------------------------------------------
# Grab G1's formatting from F1 before we repurpose F1's text, so the new
# column inherits the same header style instead of minting a new one.
$ws.Range("F1").Copy($ws.Range("G1"))

$ws.Range("C1").Value = "قابلیت عبور از فضای مجاز"
$ws.Range("E1").Value = "قابلیت عبور از فضای آزاد"
$ws.Range("F1").Value = "اندازه ورود به فضای سازه"
$ws.Range("G1").Value = "قابلیت عبور از فضای سازه"

# --- Row formatting ------------------------------------------------------
# Rows 3-6 are new legs of the same route; copy row 2's look (font/fill/
# border/number format/alignment) across A:E for each of them.
$ws.Range("A2:E2").Copy($ws.Range("A3:E3"))
$ws.Range("A2:E2").Copy($ws.Range("A4:E4"))
$ws.Range("A2:E2").Copy($ws.Range("A5:E5"))
$ws.Range("A2:E2").Copy($ws.Range("A6:E6"))
# Row 6 also uses the new F/G columns -- seed their formats from the
# existing numeric (D2) and text (C2) styled cells.
$ws.Range("D2").Copy($ws.Range("F6"))
$ws.Range("C2").Copy($ws.Range("G6"))

# Row 2's leg no longer has a "غیرمجاز" entry size -- drop F2 entirely
# (matches/fill/border included) rather than leaving it blank.
$ws.Range("F2").Clear()

# --- Row 2 (unchanged leg, still سرخس - بهرام) --------------------------
$ws.Range("A2").Value = "سرخس - بهرام"
$ws.Range("B2").Value = "گاباری 5.2"
$ws.Range("C2").Value = "غیر قابل عبور"
$ws.Range("D2").Value = 35.35533905932738
$ws.Range("E2").Value = "قابل عبور"

# --- Row 3: بهرام - باغ یک -------------------------------------------
$ws.Range("A3").Value = "بهرام - باغ یک"
$ws.Range("B3").Value = "گاباری 5.4"
$ws.Range("C3").Value = "غیر قابل عبور"
$ws.Range("D3").Value = 25.0
$ws.Range("E3").Value = "قابل عبور"

# --- Row 4: باغ یک - سواریان -----------------------------------------
$ws.Range("A4").Value = "باغ یک - سواریان"
$ws.Range("B4").Value = "گاباری 5.2"
$ws.Range("C4").Value = "غیر قابل عبور"
$ws.Range("D4").Value = 35.35533905932738
$ws.Range("E4").Value = "قابل عبور"

# --- Row 5: سواریان - نورآباد -----------------------------------------
$ws.Range("A5").Value = "سواریان - نورآباد"
$ws.Range("B5").Value = "گاباری 5.4"
$ws.Range("C5").Value = "غیر قابل عبور"
$ws.Range("D5").Value = 25.0
$ws.Range("E5").Value = "قابل عبور"

# --- Row 6: نورآباد - بندرامام خمینی (final leg, uses F/G too) --------
$ws.Range("A6").Value = "نورآباد - بندرامام خمینی"
$ws.Range("B6").Value = "گاباری ۴.۷"
$ws.Range("C6").Value = "غیر قابل عبور"
$ws.Range("D6").Value = 335.4101966249685
$ws.Range("E6").Value = "غیر قابل عبور"
$ws.Range("F6").Value = 100.0
$ws.Range("G6").Value = "غیر قابل عبور"

# New 7th column gets the same 21-char width as the others.
$ws.Columns.Item(7).ColumnWidth = 20.166666666666668
